$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990207675/1749682509529/zqlwK29jaz-camera-video-segment-driver-1749682514529.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEFZJDGT2%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T220007Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIGopD1E9f7IlBuX6xE9iW2HIADFUk%2BJajdIWXtbl3FFfAiEArLoX2FsZt%2Ftk7Ntmiq9XyAQRzBdc7oW5mKAEnUbeXWoq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCXzrifCGSZMKNdwsCq6A805aDtrhOzNwVjKeX2ofxF3r6YWfvN3q2NAfikdTa8%2BcZssD57MNYSkK4C1Vri6yRekypVXD8ESInH%2FZEepeUbRg5IWCCEW3BEJTwzId%2BMMetHahwksUhqdrCuQ1fbomFO%2FE0856Z%2FcaUXK5wHb8pVneQfpXgTc7iqTJsS2fOfhjywKZiUGmj8xqs1jsbZ4heiBG6iPqMbPB%2Bgobpt3%2FFO%2B9McjLRkPK%2FUYrc5mLEXF3bJg1fai5bpQExk1PPiyfDhEsWNIlK3mspx37BRhA7R9pOmor%2F4yCBnodpKyfAlrSHtXuaWEksutPgFrpTnaDOsau1DZbwD%2Bn3ohwdJn4EPUHDpXB9K%2FRt%2BYGB6CuGoPtbQ1r5KRz%2F1hjzkUgJ5yJfayJ9YE9ZsuS%2FE9A6Ti%2B%2BoIt5fepkS%2BOOOzUPsy2OzUsOrtEiolwELwZ8bcl9mZlKCDGKW9hVufs8FZ8BCixmSu0yZ8ICZc%2FlVnfT%2FhNj99y7lKJnv7Ckzn3xNkvNR%2B6QJWmIOuB7MXrhxurHRgkDYRF8mvDvGlET9qtmQLH2v1UWFhprbyrrKflobyPotsCD7545oSkjtQiLUw5uSswgY6pQEuduLXwj9asQfIwqqx%2FBnQjaW8WC8EpXMjXMLmPkPYn1Zf4%2B2%2Fza9sLvwQCRB0KsgEjjEP5sv6d8oxnlUqZGQDyS67HS%2BNUJRUjznTMEhIJ5W3VaeRTnwm1AtjWSyc8wjKx4RJKR4wmI79DXAGn1FqHX1CHVBcLtVLBa13PqlK%2BKzhQbI9auAQNUoGz5g2lirNdtGnTNZHo0g3VMZjx41XNuvkLD0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2006%3A00%3A07%20GMT&X-Amz-Signature=349f66a058cc77f0ba1d5869ab64d7488303f2f955077124d98c51f63b47dbf5'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990207675/1749682509529/JfEAwSUnRy-camera-video-segment-1749682514529.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEFZJDGT2%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T220006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIGopD1E9f7IlBuX6xE9iW2HIADFUk%2BJajdIWXtbl3FFfAiEArLoX2FsZt%2Ftk7Ntmiq9XyAQRzBdc7oW5mKAEnUbeXWoq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCXzrifCGSZMKNdwsCq6A805aDtrhOzNwVjKeX2ofxF3r6YWfvN3q2NAfikdTa8%2BcZssD57MNYSkK4C1Vri6yRekypVXD8ESInH%2FZEepeUbRg5IWCCEW3BEJTwzId%2BMMetHahwksUhqdrCuQ1fbomFO%2FE0856Z%2FcaUXK5wHb8pVneQfpXgTc7iqTJsS2fOfhjywKZiUGmj8xqs1jsbZ4heiBG6iPqMbPB%2Bgobpt3%2FFO%2B9McjLRkPK%2FUYrc5mLEXF3bJg1fai5bpQExk1PPiyfDhEsWNIlK3mspx37BRhA7R9pOmor%2F4yCBnodpKyfAlrSHtXuaWEksutPgFrpTnaDOsau1DZbwD%2Bn3ohwdJn4EPUHDpXB9K%2FRt%2BYGB6CuGoPtbQ1r5KRz%2F1hjzkUgJ5yJfayJ9YE9ZsuS%2FE9A6Ti%2B%2BoIt5fepkS%2BOOOzUPsy2OzUsOrtEiolwELwZ8bcl9mZlKCDGKW9hVufs8FZ8BCixmSu0yZ8ICZc%2FlVnfT%2FhNj99y7lKJnv7Ckzn3xNkvNR%2B6QJWmIOuB7MXrhxurHRgkDYRF8mvDvGlET9qtmQLH2v1UWFhprbyrrKflobyPotsCD7545oSkjtQiLUw5uSswgY6pQEuduLXwj9asQfIwqqx%2FBnQjaW8WC8EpXMjXMLmPkPYn1Zf4%2B2%2Fza9sLvwQCRB0KsgEjjEP5sv6d8oxnlUqZGQDyS67HS%2BNUJRUjznTMEhIJ5W3VaeRTnwm1AtjWSyc8wjKx4RJKR4wmI79DXAGn1FqHX1CHVBcLtVLBa13PqlK%2BKzhQbI9auAQNUoGz5g2lirNdtGnTNZHo0g3VMZjx41XNuvkLD0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2006%3A00%3A06%20GMT&X-Amz-Signature=bf11720731e1ea1f16b21d00845f9e30c6dde195d310db23c3157940654367a5'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991154589/1749681707753/UtIzP4FPML-camera-video-segment-driver-1749681712753.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEFZJDGT2%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T220006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIGopD1E9f7IlBuX6xE9iW2HIADFUk%2BJajdIWXtbl3FFfAiEArLoX2FsZt%2Ftk7Ntmiq9XyAQRzBdc7oW5mKAEnUbeXWoq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCXzrifCGSZMKNdwsCq6A805aDtrhOzNwVjKeX2ofxF3r6YWfvN3q2NAfikdTa8%2BcZssD57MNYSkK4C1Vri6yRekypVXD8ESInH%2FZEepeUbRg5IWCCEW3BEJTwzId%2BMMetHahwksUhqdrCuQ1fbomFO%2FE0856Z%2FcaUXK5wHb8pVneQfpXgTc7iqTJsS2fOfhjywKZiUGmj8xqs1jsbZ4heiBG6iPqMbPB%2Bgobpt3%2FFO%2B9McjLRkPK%2FUYrc5mLEXF3bJg1fai5bpQExk1PPiyfDhEsWNIlK3mspx37BRhA7R9pOmor%2F4yCBnodpKyfAlrSHtXuaWEksutPgFrpTnaDOsau1DZbwD%2Bn3ohwdJn4EPUHDpXB9K%2FRt%2BYGB6CuGoPtbQ1r5KRz%2F1hjzkUgJ5yJfayJ9YE9ZsuS%2FE9A6Ti%2B%2BoIt5fepkS%2BOOOzUPsy2OzUsOrtEiolwELwZ8bcl9mZlKCDGKW9hVufs8FZ8BCixmSu0yZ8ICZc%2FlVnfT%2FhNj99y7lKJnv7Ckzn3xNkvNR%2B6QJWmIOuB7MXrhxurHRgkDYRF8mvDvGlET9qtmQLH2v1UWFhprbyrrKflobyPotsCD7545oSkjtQiLUw5uSswgY6pQEuduLXwj9asQfIwqqx%2FBnQjaW8WC8EpXMjXMLmPkPYn1Zf4%2B2%2Fza9sLvwQCRB0KsgEjjEP5sv6d8oxnlUqZGQDyS67HS%2BNUJRUjznTMEhIJ5W3VaeRTnwm1AtjWSyc8wjKx4RJKR4wmI79DXAGn1FqHX1CHVBcLtVLBa13PqlK%2BKzhQbI9auAQNUoGz5g2lirNdtGnTNZHo0g3VMZjx41XNuvkLD0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2006%3A00%3A06%20GMT&X-Amz-Signature=8b97b6fda4ecdb873db00b17c5383119373bfa0925f99ad6c651761e5bb4a837'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991154589/1749681707753/418ikFdSWa-camera-video-segment-1749681712753.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEFZJDGT2%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T220006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIGopD1E9f7IlBuX6xE9iW2HIADFUk%2BJajdIWXtbl3FFfAiEArLoX2FsZt%2Ftk7Ntmiq9XyAQRzBdc7oW5mKAEnUbeXWoq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCXzrifCGSZMKNdwsCq6A805aDtrhOzNwVjKeX2ofxF3r6YWfvN3q2NAfikdTa8%2BcZssD57MNYSkK4C1Vri6yRekypVXD8ESInH%2FZEepeUbRg5IWCCEW3BEJTwzId%2BMMetHahwksUhqdrCuQ1fbomFO%2FE0856Z%2FcaUXK5wHb8pVneQfpXgTc7iqTJsS2fOfhjywKZiUGmj8xqs1jsbZ4heiBG6iPqMbPB%2Bgobpt3%2FFO%2B9McjLRkPK%2FUYrc5mLEXF3bJg1fai5bpQExk1PPiyfDhEsWNIlK3mspx37BRhA7R9pOmor%2F4yCBnodpKyfAlrSHtXuaWEksutPgFrpTnaDOsau1DZbwD%2Bn3ohwdJn4EPUHDpXB9K%2FRt%2BYGB6CuGoPtbQ1r5KRz%2F1hjzkUgJ5yJfayJ9YE9ZsuS%2FE9A6Ti%2B%2BoIt5fepkS%2BOOOzUPsy2OzUsOrtEiolwELwZ8bcl9mZlKCDGKW9hVufs8FZ8BCixmSu0yZ8ICZc%2FlVnfT%2FhNj99y7lKJnv7Ckzn3xNkvNR%2B6QJWmIOuB7MXrhxurHRgkDYRF8mvDvGlET9qtmQLH2v1UWFhprbyrrKflobyPotsCD7545oSkjtQiLUw5uSswgY6pQEuduLXwj9asQfIwqqx%2FBnQjaW8WC8EpXMjXMLmPkPYn1Zf4%2B2%2Fza9sLvwQCRB0KsgEjjEP5sv6d8oxnlUqZGQDyS67HS%2BNUJRUjznTMEhIJ5W3VaeRTnwm1AtjWSyc8wjKx4RJKR4wmI79DXAGn1FqHX1CHVBcLtVLBa13PqlK%2BKzhQbI9auAQNUoGz5g2lirNdtGnTNZHo0g3VMZjx41XNuvkLD0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2006%3A00%3A06%20GMT&X-Amz-Signature=1bef96c1fdc2b528dbc4a12e038b2b0e8e4de8d0f7c582a8646d12aa2039c091'

$wb.Save()